$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.203.15"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "'2.918.31"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'569.97"
$ws.Range("E5").Value = "  -3.08%  "
$ws.Range("D6").Value = "'144.39"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'2.917.77"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "'6.96"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "'0.149"
$ws.Range("E11").Value = "  -2.30%  "
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "'0.0000240"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "'32.69"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "'3.401.52"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'62.117.48"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "'2.915.08"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("D20").Value = "'432.10"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "'13.10"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").Value = "'0.655"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Value = "'78.94"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("D26").Value = "'10.23"
$ws.Range("E26").Value = "  -2.70%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("E29").Value = "  +6.11%  "
$ws.Range("D30").Value = "'6.98"
$ws.Range("E30").Value = "  -4.95%  "
$ws.Range("D31").Value = "'2.52"
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("D35").Value = "'25.75"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "'0.956"
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("D38").Value = "'48.89"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").Value = "'2.93"
$ws.Range("E39").Value = "  -6.03%  "
$ws.Range("E40").Value = "  -4.83%  "
$ws.Range("D41").Value = "'41.21"
$ws.Range("E41").Value = "  +4.75%  "
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("D45").Value = "'2.728.00"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("D46").Value = "'133.56"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").Value = "'349.10"
$ws.Range("E48").Value = "  +0.45%  "

# Rows 50 and 51 swap content (Stellar <-> FLOKI) with updated values
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000216"
$ws.Range("E50").Value = "  +10.59%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.103"
$ws.Range("E51").Value = "  -0.94%  "
